{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per diff):\n//  - para 0: date 23.11.24 -> 22.11.24\n//  - para 1: title replaced\n//  - para 2: body replaced (new trailing-space preserved text)\n//  - para 3 (Heading2 \"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\"): removed\n//  - para 4: body replaced\n//  - para 5 (Heading2 \"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\"): removed\n//  - para 6: body replaced\n//  - para 7: body replaced\n//  - para 8: body replaced\n//  - paras 9..27 (19 paragraphs): removed entirely\n//  - para 28 (last, the arxiv URL): text replaced\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// New text for the paragraphs that are kept (matched by their ORIGINAL text).\nconst replacements = [\n  {\n    find: \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.11.24: \u26a1\ufe0f\ud83d\ude80\",\n    text: \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -22.11.24: \u26a1\ufe0f\ud83d\ude80\",\n  },\n  {\n    find: \"Table Meets LLM: Can Large Language Models Understand Structured Table Data? A Benchmark and Empirical Study\",\n    text: \"The Unreasonable Ineffectiveness of the Deeper Layers\",\n  },\n  {\n    find: \"\u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05e1\u05d5\u05e7\u05e8 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e0\u05d2\u05e2\u05ea \u05d1\u05d5(\u05d1\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea) \u05d5\u05d4\u05d5\u05d0 \u05d3\u05d0\u05d8\u05d4 \u05d8\u05d1\u05dc\u05d0\u05d9. \u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05e9\u05d0\u05dc\u05d4 \u05de\u05e8\u05ea\u05e7\u05ea - \u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (LLMs) \u05db\u05de\u05d5 GPT \u05d1\u05d0\u05de\u05ea \u05de\u05d1\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4 \u05d1\u05d8\u05d1\u05dc\u05d0\u05d5\u05ea?\",\n    text: \"\u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d9\u05dc \u05e9\u05dc\u05d0 \u05d9\u05e7\u05e9\u05d4 \u05e2\u05dc\u05d9\u05db\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d3\u05d9 \u05d1\u05e1\u05d5\u05e4\u05f4\u05e9. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05de\u05d0\u05d5\u05d3 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05dc\u05e7\u05e6\u05e5 \u05e9\u05db\u05d1\u05d5\u05ea \u05d1\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d9\u05dd \u05e2\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e9\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e9\u05dc\u05e0\u05d5 \u05d5\u05d2\u05dd \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d1\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05d7\u05e8\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d9\u05dd \u05e2\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e9\u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de\u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e9\u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05de\u05e0\u05d2\u05e0\u05d5\u05df attention \u05d5\u05e9\u05ea\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea feed-forward (\u05d4\u05e9\u05e0\u05d9\u05d9\u05d4 \u05de\u05d4\u05df \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea). \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d9\u05e9 \u05e9\u05db\u05d1\u05d5\u05ea \u05e0\u05e8\u05de\u05d5\u05dc \u05d5\u05d7\u05d9\u05d1\u05d5\u05e8\u05d9 residual (\u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05de\u05d7\u05d5\u05d1\u05e8 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea). \",\n  },\n  {\n    find: \"\u05d1\u05e9\u05e0\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea, LLMs \u05d4\u05e4\u05db\u05d5 \u05dc\u05db\u05dc\u05d9 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea. \u05d0\u05d1\u05dc \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d4\u05dd \u05de\u05e6\u05d5\u05d9\u05e0\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc) \u05d1\u05d4\u05d1\u05e0\u05ea \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8), \u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05e2\u05d3\u05d9\u05d9\u05df \u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d4 \u05dc\u05e2\u05d5\u05de\u05e7 \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8\",\n    text: \"\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05de\u05db\u05d9\u05dc\u05d9\u05dd \u05e2\u05e9\u05e8\u05d5\u05ea \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05d1\u05dc\u05d5\u05e7\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e9\u05db\u05de\u05d5\u05d1\u05df \u05de\u05e9\u05dc\u05d9\u05da \u05e2\u05dc \u05db\u05de\u05d5\u05ea \u05d4\u05d6\u05de\u05df \u05d5\u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd \u05d4\u05e0\u05d3\u05e8\u05e9\u05d9\u05dd \u05dc\u05d4\u05e4\u05e2\u05dc\u05ea\u05dd, \u05d1\u05e2\u05d9\u05e7\u05e8 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d2\u05e0\u05e8\u05d5\u05d8. \u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05dc\u05e7\u05e6\u05e5 \u05db\u05de\u05d4 \u05d1\u05dc\u05d5\u05e7\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd \u05e9\u05db\u05de\u05d5\u05d1\u05df \u05d9\u05e7\u05d8\u05d9\u05df \u05d0\u05ea \u05d6\u05de\u05df \u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05e0\u05d3\u05e8\u05e9 \u05dc\u05d9\u05e6\u05d9\u05e8\u05d4 \u05d4\u05e4\u05dc\u05d8. \u05d0\u05d1\u05dc \u05d0\u05d9\u05d6\u05d4 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05dc\u05d1\u05d7\u05d5\u05e8 \u05db\u05da \u05e9\u05d4\u05e4\u05d2\u05d9\u05e2\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d4\u05d9\u05d4 \u05de\u05d9\u05e0\u05d9\u05de\u05dc\u05d9\u05ea.\",\n  },\n  {\n    find: \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05de\u05d3\u05d3 \u05d7\u05d3\u05e9 \u05e9\u05e0\u05e7\u05e8\u05d0 (SUC (Structural Understanding Capabilities \u05e9\u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea. \u05d4\u05de\u05d3\u05d3 \u05db\u05d5\u05dc\u05dc \u05e9\u05d1\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea:\",\n    text: \"\u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d4\u05d2\u05e8\u05e3 \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9 \u05e9\u05dc \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05dc\u05d0 \u05de\u05e2\u05d8 \u05d7\u05d9\u05d1\u05d5\u05e8\u05d9 residual \u05d8\u05d1\u05e2\u05d9 \u05dc\u05d1\u05d7\u05d5\u05e8 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd \u05e9\u05dc\u05d0 \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05d4\u05e8\u05d1\u05d4 \u05dc\u05e4\u05dc\u05d8 \u05d4\u05d1\u05dc\u05d5\u05e7 \u05d4\u05e0\u05de\u05e6\u05d0 \u05dc\u05e4\u05e0\u05d9\u05d4\u05dd \u05d1\u05de\u05d5\u05d3\u05dc. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05d3\u05dc\u05ea\u05d0 \u05e9\u05e0\u05d5\u05ea\u05e0\u05d9\u05dd \u05d4\u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 \u05d6\u05e0\u05d9\u05d7\u05d4 \u05d0\u05d6 \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05e2\u05d9\u05e3 \u05d0\u05d5\u05ea\u05dd \u05d1\u05dc\u05d9 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05e8\u05e6\u05d9\u05e0\u05d9\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd.\",\n  },\n  {\n    find: \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4\",\n    text: \"\u05d4\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05e0\u05d9\u05ea\u05df \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d6\u05d4? \u05d4\u05d0\u05de\u05ea \u05d9\u05e9 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d6\u05d4 \u05d5\u05de\u05d0\u05de\u05e8 \u05d1\u05d7\u05e8 \u05dc\u05d4\u05e9\u05d5\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d1\u05dc\u05d5\u05e7 l \u05e2\u05dd \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d1\u05dc\u05d5\u05e7 l+n (\u05d0\u05e0\u05d5 \u05de\u05d5\u05d7\u05e7\u05d9\u05dd n \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd) \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05de\u05d5\u05d3\u05d9\u05e4\u05d9\u05e7\u05e6\u05d9\u05d4 \u05e7\u05d8\u05e0\u05d4 \u05e9\u05dc \u05de\u05e8\u05d7\u05e7 \u05e7\u05d5\u05e1\u05d9\u05d9\u05df (\u05d4\u05d7\u05dc\u05d9\u05e4\u05d5 cos \u05d1-arccos \u05d5\u05d7\u05d9\u05dc\u05e7\u05d5 \u05d1-pi \u05db\u05d3\u05d9 \u05dc\u05d2\u05e8\u05d5\u05dd \u05dc\u05de\u05d3\u05d3 \u05d4\u05d6\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05d1\u05d9\u05df 0 \u05dc 1). \u05d1\u05d0\u05d5\u05e4\u05df \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9 n \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e2\u05dd \u05d3\u05de\u05d9\u05d5\u05df \u05d2\u05d1\u05d5\u05d4 \u05de\u05d0\u05d5\u05d3 \u05dc\u05d1\u05dc\u05d5\u05e7 \u05e9\u05e7\u05d5\u05d3\u05dd \u05dc\u05d4\u05dd (\u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05e4\u05dc\u05d8) \u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 \u05de\u05d5\u05e2\u05de\u05d3\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05dc\u05e7\u05d9\u05e6\u05d5\u05e5 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d1\u05dc\u05d5\u05e7 \u05d4\u05ea\u05d7\u05dc\u05ea\u05d9 l \u05d5\u05de\u05e1\u05e4\u05e8 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05dc\u05e7\u05d9\u05e6\u05d5\u05e5 n \u05e2\u05dd \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05d4\u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8). \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05de\u05d7\u05d5\u05e9\u05d1 \u05e2\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05e2\u05d1\u05d5\u05e8 \u05db\u05de\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d4.\",\n  },\n  {\n    find: \"\u05d0\u05d9\u05ea\u05d5\u05e8 \u05ea\u05d0\u05d9\u05dd \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd\",\n    text: \"\u05dc\u05d0\u05d7\u05e8 \u05d4\u05de\u05d7\u05d9\u05e7\u05d4 \u05e0\u05d9\u05ea\u05df \u05dc\u05e2\u05e9\u05d5\u05ea \u05dc\u05de\u05d5\u05d3\u05dc \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e7\u05dc\u05d9\u05dc \u05d5\u05dc\u05d8\u05e2\u05e0\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05de\u05d7\u05d5\u05e7 \u05db\u05db\u05d4 \u05e2\u05dc \u05d7\u05e6\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea \u05d8\u05e0\u05e8\u05e1\u05e4\u05d5\u05e8\u05de\u05d9\u05dd (\u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4) \u05d1\u05dc\u05d9 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05e8\u05e6\u05d9\u05e0\u05d9\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd).\",\n  },\n  {\n    find: \"https://arxiv.org/abs/2305.13062\",\n    text: \"https://arxiv.org/abs/2403.17887\",\n  },\n];\n\n// Paragraphs whose entire text exactly matches one of these should be removed.\nconst toDelete = new Set([\n  \"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\",\n  \"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\",\n  \"\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05d5\u05da (\u05de\u05d9\u05e7\u05d5\u05dd \u05dc\u05e2\u05e8\u05da)\",\n  \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05de\u05d5\u05d3\u05d5\u05ea\",\n  \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d5\u05ea\",\n  \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4\",\n  \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd\",\n  \"\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05d0\u05ea GPT-3.5 \u05d5-GPT-4 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05ea\u05d5\u05da \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05dc\u05d8 (HTML, JSON, CSV \u05d5\u05e2\u05d5\u05d3).\",\n  \"\u05de\u05d4 \u05d4\u05dd \u05d2\u05d9\u05dc\u05d5?\",\n  \"\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d5\u05ea! \u05d4\u05e0\u05d4 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:\",\n  \"HTML \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e4\u05d5\u05e8\u05de\u05d8 \u05f4\u05d4\u05e0\u05d5\u05d7\u05f4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05d4\u05e6\u05d2\u05ea \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05dc-LLMs\",\n  \"\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d9\u05d7\u05e1\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4, \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd) \u05d0\u05da \u05e0\u05db\u05e9\u05dc\u05d5 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4, \u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d4 \u05e4\u05e9\u05d5\u05d8, \u05d7\u05d9\u05e4\u05d5\u05e9 \u05ea\u05d0 \u05d1\u05d5\u05d3\u05d3)\",\n  \"\u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d4\u05e9\u05ea\u05e4\u05e8\u05d5 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05e2\u05dd \u05d3\u05d5\u05d2\u05de\u05d4 \u05d0\u05d7\u05ea (one-shot) \u05dc\u05e2\u05d5\u05de\u05ea \u05d0\u05e4\u05e1 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea\",\n  \"\u05d4\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9: Self-augmented Prompting\",\n  \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0\u05ea \\\"self-augmented prompting\\\" \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05e7\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d7\u05d9\u05dc\u05d4 \u05dc\u05d6\u05d4\u05d5\u05ea \u05de\u05d9\u05d3\u05e2 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05d1\u05d8\u05d1\u05dc\u05d4 (\u05db\u05de\u05d5 \u05d8\u05d5\u05d5\u05d7\u05d9 \u05e2\u05e8\u05db\u05d9\u05dd) \u05d5\u05d0\u05d6 \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05de\u05d9\u05d3\u05e2 \u05d4\u05d6\u05d4 \u05db\u05d3\u05d9 \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d1\u05e0\u05e6'\u05de\u05d0\u05e8\u05e7\u05d9\u05dd)\",\n  \"\u05e1\u05d9\u05db\u05d5\u05dd:\",\n  \"\u05d0\u05e0\u05d9 \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05e8\u05ea\u05e7. \u05d4\u05d5\u05d0 \u05de\u05e8\u05d0\u05d4 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05d4\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d4\u05e2\u05e6\u05d5\u05de\u05d4 \u05d1-LLMs, \u05d9\u05e9 \u05e2\u05d3\u05d9\u05d9\u05df \u05e4\u05e2\u05e8\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc\u05d4\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4. \u05d6\u05d4 \u05de\u05d6\u05db\u05d9\u05e8 \u05dc\u05e0\u05d5 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05e8\u05e9\u05d9\u05de\u05d9\u05dd, \u05d4\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e8\u05d7\u05d5\u05e7\u05d9\u05dd \u05de\u05d4\u05d1\u05e0\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05de\u05d1\u05e0\u05d9\u05dd \u05d5\u05d9\u05d7\u05e1\u05d9\u05dd \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4.\",\n  \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5 \u05e2\u05d1\u05d5\u05d3\u05d4 \u05dc\u05d0 \u05e8\u05e2\u05d4 \u05d1\u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d3\u05d3\u05d9\u05dd \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05d9\u05e2\u05d6\u05e8\u05d5 \u05dc\u05e7\u05d4\u05d9\u05dc\u05d4 \u05dc\u05d4\u05de\u05e9\u05d9\u05da \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4 \u05e9\u05dc\u05d4\u05dd \u05dc-prompting \u05d4\u05d9\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05d1\u05dc \u05d0\u05e4\u05e7\u05d8\u05d9\u05d1\u05d9\u05ea, \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd - \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e4\u05e8\u05e7\u05d8\u05d9\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d9\u05d9\u05e9\u05dd \u05de\u05d9\u05d3.\",\n  \"\u05de\u05d9\u05dc\u05d4 \u05d0\u05d7\u05e8\u05d5\u05e0\u05d4\",\n  \"\u05d0\u05dd \u05d0\u05ea\u05dd \u05e2\u05d5\u05d1\u05d3\u05d9\u05dd \u05e2\u05dd \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05d5-LLMs, \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d7\u05d5\u05d1\u05d4. \u05d4\u05d5\u05d0 \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05d5\u05db\u05dc\u05d9\u05dd \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05d9\u05dd. \u05d4\u05e7\u05d5\u05d3 \u05d5\u05d4\u05d3\u05d0\u05d8\u05d4 \u05d6\u05de\u05d9\u05e0\u05d9\u05dd \u05d1-GitHub, \u05d0\u05d6 \u05d0\u05ea\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05ea\u05d7\u05d9\u05dc \u05dc\u05e9\u05d7\u05e7 \u05e2\u05dd \u05d6\u05d4 \u05d9\u05e9\u05e8.\",\n  \"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05d9\u05d4\u05d9\u05d4 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05d9\u05da \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05d9\u05e9\u05e4\u05d9\u05e2\u05d5 \u05e2\u05dc \u05d4\u05d3\u05d5\u05e8 \u05d4\u05d1\u05d0 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05d0\u05dd \u05e0\u05e8\u05d0\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05de\u05ea\u05d5\u05db\u05e0\u05e0\u05d9\u05dd \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d4\u05d1\u05e0\u05ea \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4?\",\n]);\n\n// First pass: apply text replacements (paragraph identity doesn't change on text edit).\nfor (const item of items) {\n  const current = item.text;\n  const match = replacements.find((r) => current === r.find);\n  if (match) {\n    item.insertText(match.text, \"Replace\");\n  }\n}\n\n// Second pass: delete paragraphs that must disappear entirely. Walk from the\n// end backwards so deleting doesn't disturb indices of paragraphs not yet\n// processed.\nfor (let i = items.length - 1; i >= 0; i--) {\n  const current = items[i].text;\n  if (toDelete.has(current)) {\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per diff):\n#  - para 0: date 23.11.24 -> 22.11.24\n#  - para 1: title replaced\n#  - para 2: body replaced (new trailing-space preserved text)\n#  - para 3 (Heading2 \"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\"): removed\n#  - para 4: body replaced\n#  - para 5 (Heading2 \"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\"): removed\n#  - para 6: body replaced\n#  - para 7: body replaced\n#  - para 8: body replaced\n#  - paras 9..27 (19 paragraphs): removed entirely\n#  - para 28 (last, the arxiv URL): text replaced\n\n$d = $word.ActiveDocument\n\n# New text for the paragraphs that are kept, keyed by their ORIGINAL text.\n$replacements = @{\n    \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.11.24: \u26a1\ufe0f\ud83d\ude80\" = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -22.11.24: \u26a1\ufe0f\ud83d\ude80\";\n    \"Table Meets LLM: Can Large Language Models Understand Structured Table Data? A Benchmark and Empirical Study\" = \"The Unreasonable Ineffectiveness of the Deeper Layers\";\n    \"\u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05e1\u05d5\u05e7\u05e8 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e0\u05d2\u05e2\u05ea \u05d1\u05d5(\u05d1\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea) \u05d5\u05d4\u05d5\u05d0 \u05d3\u05d0\u05d8\u05d4 \u05d8\u05d1\u05dc\u05d0\u05d9. \u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05e9\u05d0\u05dc\u05d4 \u05de\u05e8\u05ea\u05e7\u05ea - \u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (LLMs) \u05db\u05de\u05d5 GPT \u05d1\u05d0\u05de\u05ea \u05de\u05d1\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4 \u05d1\u05d8\u05d1\u05dc\u05d0\u05d5\u05ea?\" = \"\u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d9\u05dc \u05e9\u05dc\u05d0 \u05d9\u05e7\u05e9\u05d4 \u05e2\u05dc\u05d9\u05db\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d3\u05d9 \u05d1\u05e1\u05d5\u05e4\u05f4\u05e9. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05de\u05d0\u05d5\u05d3 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05dc\u05e7\u05e6\u05e5 \u05e9\u05db\u05d1\u05d5\u05ea \u05d1\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05d1\u05d5\u05e1\u05e1\u05d9\u05dd \u05e2\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e9\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e9\u05dc\u05e0\u05d5 \u05d5\u05d2\u05dd \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d1\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05d7\u05e8\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d9\u05dd \u05e2\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e9\u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de\u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e9\u05db\u05dc \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05de\u05e0\u05d2\u05e0\u05d5\u05df attention \u05d5\u05e9\u05ea\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea feed-forward (\u05d4\u05e9\u05e0\u05d9\u05d9\u05d4 \u05de\u05d4\u05df \u05dc\u05d9\u05e0\u05d0\u05e8\u05d9\u05ea). \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d9\u05e9 \u05e9\u05db\u05d1\u05d5\u05ea \u05e0\u05e8\u05de\u05d5\u05dc \u05d5\u05d7\u05d9\u05d1\u05d5\u05e8\u05d9 residual (\u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05db\u05dc \u05e9\u05db\u05d1\u05d4 \u05de\u05d7\u05d5\u05d1\u05e8 \u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea). \";\n    \"\u05d1\u05e9\u05e0\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea, LLMs \u05d4\u05e4\u05db\u05d5 \u05dc\u05db\u05dc\u05d9 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea. \u05d0\u05d1\u05dc \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d4\u05dd \u05de\u05e6\u05d5\u05d9\u05e0\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc) \u05d1\u05d4\u05d1\u05e0\u05ea \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8), \u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05e2\u05d3\u05d9\u05d9\u05df \u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d4 \u05dc\u05e2\u05d5\u05de\u05e7 \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8\" = \"\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05de\u05db\u05d9\u05dc\u05d9\u05dd \u05e2\u05e9\u05e8\u05d5\u05ea \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05d1\u05dc\u05d5\u05e7\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e9\u05db\u05de\u05d5\u05d1\u05df \u05de\u05e9\u05dc\u05d9\u05da \u05e2\u05dc \u05db\u05de\u05d5\u05ea \u05d4\u05d6\u05de\u05df \u05d5\u05d4\u05de\u05e9\u05d0\u05d1\u05d9\u05dd \u05d4\u05e0\u05d3\u05e8\u05e9\u05d9\u05dd \u05dc\u05d4\u05e4\u05e2\u05dc\u05ea\u05dd, \u05d1\u05e2\u05d9\u05e7\u05e8 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d2\u05e0\u05e8\u05d5\u05d8. \u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05de\u05e6\u05d9\u05e2 \u05d3\u05e8\u05da \u05dc\u05e7\u05e6\u05e5 \u05db\u05de\u05d4 \u05d1\u05dc\u05d5\u05e7\u05d9 \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd \u05e9\u05db\u05de\u05d5\u05d1\u05df \u05d9\u05e7\u05d8\u05d9\u05df \u05d0\u05ea \u05d6\u05de\u05df \u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05e0\u05d3\u05e8\u05e9 \u05dc\u05d9\u05e6\u05d9\u05e8\u05d4 \u05d4\u05e4\u05dc\u05d8. \u05d0\u05d1\u05dc \u05d0\u05d9\u05d6\u05d4 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05dc\u05d1\u05d7\u05d5\u05e8 \u05db\u05da \u05e9\u05d4\u05e4\u05d2\u05d9\u05e2\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d4\u05d9\u05d4 \u05de\u05d9\u05e0\u05d9\u05de\u05dc\u05d9\u05ea.\";\n    \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05de\u05d3\u05d3 \u05d7\u05d3\u05e9 \u05e9\u05e0\u05e7\u05e8\u05d0 (SUC (Structural Understanding Capabilities \u05e9\u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea. \u05d4\u05de\u05d3\u05d3 \u05db\u05d5\u05dc\u05dc \u05e9\u05d1\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea:\" = \"\u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d4\u05d2\u05e8\u05e3 \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9 \u05e9\u05dc \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8 \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05dc\u05d0 \u05de\u05e2\u05d8 \u05d7\u05d9\u05d1\u05d5\u05e8\u05d9 residual \u05d8\u05d1\u05e2\u05d9 \u05dc\u05d1\u05d7\u05d5\u05e8 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd \u05e9\u05dc\u05d0 \u05de\u05d5\u05e1\u05d9\u05e4\u05d9\u05dd \u05d4\u05e8\u05d1\u05d4 \u05dc\u05e4\u05dc\u05d8 \u05d4\u05d1\u05dc\u05d5\u05e7 \u05d4\u05e0\u05de\u05e6\u05d0 \u05dc\u05e4\u05e0\u05d9\u05d4\u05dd \u05d1\u05de\u05d5\u05d3\u05dc. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05d3\u05dc\u05ea\u05d0 \u05e9\u05e0\u05d5\u05ea\u05e0\u05d9\u05dd \u05d4\u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5 \u05d6\u05e0\u05d9\u05d7\u05d4 \u05d0\u05d6 \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05e2\u05d9\u05e3 \u05d0\u05d5\u05ea\u05dd \u05d1\u05dc\u05d9 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05e8\u05e6\u05d9\u05e0\u05d9\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd.\";\n    \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4\" = \"\u05d4\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05e0\u05d9\u05ea\u05df \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d6\u05d4? \u05d4\u05d0\u05de\u05ea \u05d9\u05e9 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d6\u05d4 \u05d5\u05de\u05d0\u05de\u05e8 \u05d1\u05d7\u05e8 \u05dc\u05d4\u05e9\u05d5\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d1\u05dc\u05d5\u05e7 l \u05e2\u05dd \u05d4\u05e4\u05dc\u05d8 \u05e9\u05dc \u05d4\u05d1\u05dc\u05d5\u05e7 l+n (\u05d0\u05e0\u05d5 \u05de\u05d5\u05d7\u05e7\u05d9\u05dd n \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e8\u05e6\u05d5\u05e4\u05d9\u05dd) \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05de\u05d5\u05d3\u05d9\u05e4\u05d9\u05e7\u05e6\u05d9\u05d4 \u05e7\u05d8\u05e0\u05d4 \u05e9\u05dc \u05de\u05e8\u05d7\u05e7 \u05e7\u05d5\u05e1\u05d9\u05d9\u05df (\u05d4\u05d7\u05dc\u05d9\u05e4\u05d5 cos \u05d1-arccos \u05d5\u05d7\u05d9\u05dc\u05e7\u05d5 \u05d1-pi \u05db\u05d3\u05d9 \u05dc\u05d2\u05e8\u05d5\u05dd \u05dc\u05de\u05d3\u05d3 \u05d4\u05d6\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05d1\u05d9\u05df 0 \u05dc 1). \u05d1\u05d0\u05d5\u05e4\u05df \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9 n \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05e2\u05dd \u05d3\u05de\u05d9\u05d5\u05df \u05d2\u05d1\u05d5\u05d4 \u05de\u05d0\u05d5\u05d3 \u05dc\u05d1\u05dc\u05d5\u05e7 \u05e9\u05e7\u05d5\u05d3\u05dd \u05dc\u05d4\u05dd (\u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d4\u05e4\u05dc\u05d8) \u05e0\u05d1\u05d7\u05e8\u05d9\u05dd \u05d1\u05ea\u05d5\u05e8 \u05de\u05d5\u05e2\u05de\u05d3\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05dc\u05e7\u05d9\u05e6\u05d5\u05e5 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d5\u05d7\u05e8\u05d9\u05dd \u05d1\u05dc\u05d5\u05e7 \u05d4\u05ea\u05d7\u05dc\u05ea\u05d9 l \u05d5\u05de\u05e1\u05e4\u05e8 \u05d1\u05dc\u05d5\u05e7\u05d9\u05dd \u05dc\u05e7\u05d9\u05e6\u05d5\u05e5 n \u05e2\u05dd \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05d4\u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8). \u05d4\u05d3\u05de\u05d9\u05d5\u05df \u05de\u05d7\u05d5\u05e9\u05d1 \u05e2\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05d8\u05d5\u05e7\u05df \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05e2\u05d1\u05d5\u05e8 \u05db\u05de\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d4.\";\n    \"\u05d0\u05d9\u05ea\u05d5\u05e8 \u05ea\u05d0\u05d9\u05dd \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd\" = \"\u05dc\u05d0\u05d7\u05e8 \u05d4\u05de\u05d7\u05d9\u05e7\u05d4 \u05e0\u05d9\u05ea\u05df \u05dc\u05e2\u05e9\u05d5\u05ea \u05dc\u05de\u05d5\u05d3\u05dc \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e7\u05dc\u05d9\u05dc \u05d5\u05dc\u05d8\u05e2\u05e0\u05ea \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05e0\u05d9\u05ea\u05df \u05dc\u05de\u05d7\u05d5\u05e7 \u05db\u05db\u05d4 \u05e2\u05dc \u05d7\u05e6\u05d9 \u05e9\u05db\u05d1\u05d5\u05ea \u05d8\u05e0\u05e8\u05e1\u05e4\u05d5\u05e8\u05de\u05d9\u05dd (\u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4) \u05d1\u05dc\u05d9 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05e8\u05e6\u05d9\u05e0\u05d9\u05ea \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd).\";\n    \"https://arxiv.org/abs/2305.13062\" = \"https://arxiv.org/abs/2403.17887\";\n}\n\n# Paragraphs whose entire text exactly matches one of these should be removed.\n$toDelete = @(\n    \"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\",\n    \"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\",\n    \"\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05d5\u05da (\u05de\u05d9\u05e7\u05d5\u05dd \u05dc\u05e2\u05e8\u05da)\",\n    \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05de\u05d5\u05d3\u05d5\u05ea\",\n    \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d5\u05ea\",\n    \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4\",\n    \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd\",\n    \"\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05d0\u05ea GPT-3.5 \u05d5-GPT-4 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05ea\u05d5\u05da \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05dc\u05d8 (HTML, JSON, CSV \u05d5\u05e2\u05d5\u05d3).\",\n    \"\u05de\u05d4 \u05d4\u05dd \u05d2\u05d9\u05dc\u05d5?\",\n    \"\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d5\u05ea! \u05d4\u05e0\u05d4 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:\",\n    \"HTML \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e4\u05d5\u05e8\u05de\u05d8 \u05f4\u05d4\u05e0\u05d5\u05d7\u05f4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05d4\u05e6\u05d2\u05ea \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05dc-LLMs\",\n    \"\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d9\u05d7\u05e1\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4, \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd) \u05d0\u05da \u05e0\u05db\u05e9\u05dc\u05d5 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4, \u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d4 \u05e4\u05e9\u05d5\u05d8, \u05d7\u05d9\u05e4\u05d5\u05e9 \u05ea\u05d0 \u05d1\u05d5\u05d3\u05d3)\",\n    \"\u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d4\u05e9\u05ea\u05e4\u05e8\u05d5 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05e2\u05dd \u05d3\u05d5\u05d2\u05de\u05d4 \u05d0\u05d7\u05ea (one-shot) \u05dc\u05e2\u05d5\u05de\u05ea \u05d0\u05e4\u05e1 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea\",\n    \"\u05d4\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9: Self-augmented Prompting\",\n    \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0\u05ea \"\"self-augmented prompting\"\" \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05e7\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d7\u05d9\u05dc\u05d4 \u05dc\u05d6\u05d4\u05d5\u05ea \u05de\u05d9\u05d3\u05e2 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05d1\u05d8\u05d1\u05dc\u05d4 (\u05db\u05de\u05d5 \u05d8\u05d5\u05d5\u05d7\u05d9 \u05e2\u05e8\u05db\u05d9\u05dd) \u05d5\u05d0\u05d6 \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05de\u05d9\u05d3\u05e2 \u05d4\u05d6\u05d4 \u05db\u05d3\u05d9 \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d1\u05e0\u05e6'\u05de\u05d0\u05e8\u05e7\u05d9\u05dd)\",\n    \"\u05e1\u05d9\u05db\u05d5\u05dd:\",\n    \"\u05d0\u05e0\u05d9 \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05e8\u05ea\u05e7. \u05d4\u05d5\u05d0 \u05de\u05e8\u05d0\u05d4 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05d4\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d4\u05e2\u05e6\u05d5\u05de\u05d4 \u05d1-LLMs, \u05d9\u05e9 \u05e2\u05d3\u05d9\u05d9\u05df \u05e4\u05e2\u05e8\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc\u05d4\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4. \u05d6\u05d4 \u05de\u05d6\u05db\u05d9\u05e8 \u05dc\u05e0\u05d5 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05e8\u05e9\u05d9\u05de\u05d9\u05dd, \u05d4\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e8\u05d7\u05d5\u05e7\u05d9\u05dd \u05de\u05d4\u05d1\u05e0\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05de\u05d1\u05e0\u05d9\u05dd \u05d5\u05d9\u05d7\u05e1\u05d9\u05dd \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4.\",\n    \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5 \u05e2\u05d1\u05d5\u05d3\u05d4 \u05dc\u05d0 \u05e8\u05e2\u05d4 \u05d1\u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d3\u05d3\u05d9\u05dd \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05d9\u05e2\u05d6\u05e8\u05d5 \u05dc\u05e7\u05d4\u05d9\u05dc\u05d4 \u05dc\u05d4\u05de\u05e9\u05d9\u05da \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4 \u05e9\u05dc\u05d4\u05dd \u05dc-prompting \u05d4\u05d9\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05d1\u05dc \u05d0\u05e4\u05e7\u05d8\u05d9\u05d1\u05d9\u05ea, \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd - \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e4\u05e8\u05e7\u05d8\u05d9\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d9\u05d9\u05e9\u05dd \u05de\u05d9\u05d3.\",\n    \"\u05de\u05d9\u05dc\u05d4 \u05d0\u05d7\u05e8\u05d5\u05e0\u05d4\",\n    \"\u05d0\u05dd \u05d0\u05ea\u05dd \u05e2\u05d5\u05d1\u05d3\u05d9\u05dd \u05e2\u05dd \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05d5-LLMs, \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d7\u05d5\u05d1\u05d4. \u05d4\u05d5\u05d0 \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05d5\u05db\u05dc\u05d9\u05dd \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05d9\u05dd. \u05d4\u05e7\u05d5\u05d3 \u05d5\u05d4\u05d3\u05d0\u05d8\u05d4 \u05d6\u05de\u05d9\u05e0\u05d9\u05dd \u05d1-GitHub, \u05d0\u05d6 \u05d0\u05ea\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05ea\u05d7\u05d9\u05dc \u05dc\u05e9\u05d7\u05e7 \u05e2\u05dd \u05d6\u05d4 \u05d9\u05e9\u05e8.\",\n    \"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05d9\u05d4\u05d9\u05d4 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05d9\u05da \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05d9\u05e9\u05e4\u05d9\u05e2\u05d5 \u05e2\u05dc \u05d4\u05d3\u05d5\u05e8 \u05d4\u05d1\u05d0 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05d0\u05dd \u05e0\u05e8\u05d0\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05de\u05ea\u05d5\u05db\u05e0\u05e0\u05d9\u05dd \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d4\u05d1\u05e0\u05ea \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4?\"\n)\n\n# Walk paragraphs from last to first so deleting never shifts the index of a\n# paragraph we haven't processed yet.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    $raw = $p.Range.Text\n    $current = $raw.TrimEnd([char]13, [char]7)\n\n    if ($toDelete -contains $current) {\n        $p.Range.Delete()\n        continue\n    }\n\n    if ($replacements.ContainsKey($current)) {\n        $p.Range.Text = $replacements[$current]\n    }\n}\n"}
